$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "Joint regime area ..." rows (36-40)
$ws.Rows("36:40").Delete()

# Update weather-cluster solar assignment values (rows 2-35)
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 0.007669616519174041
$ws.Range("G2").Value = 0.01514023330851326
$ws.Range("I2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.2109374999999996
$ws.Range("F3").Value = 0.003539823008849557
$ws.Range("G3").Value = 0.03747828245222147
$ws.Range("I3").Value = 0.2433173406442785
$ws.Range("B4").Value = 0.1335453100158984
$ws.Range("C4").Value = 0.01694915254237288
$ws.Range("D4").Value = 0.1395348837209302
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = 0.05132743362831858
$ws.Range("G5").Value = 0.02084884586746091
$ws.Range("H5").Value = 0.0009319664492078285
$ws.Range("I5").Value = 0.1459904043865669
$ws.Range("F6").Value = 0.01179941002949852
$ws.Range("G6").Value = 0.02208984859766693
$ws.Range("I6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0.1233521657250469
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.3370535714285727
$ws.Range("F8").Value = 0.001769911504424779
$ws.Range("G8").Value = 0.03772648299826267
$ws.Range("I8").Value = 0.1720356408498981
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0.08474576271186439
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 0.02123893805309734
$ws.Range("H9").Value = 0.1621621621621618
$ws.Range("B10").Value = 0.1462639109697935
$ws.Range("C10").Value = 0.04237288135593222
$ws.Range("D10").Value = 0.2325581395348837
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("F11").Value = 0.01238938053097345
$ws.Range("H11").Value = 0.004659832246039142
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0.01883239171374764
$ws.Range("D12").Value = 0
$ws.Range("F12").Value = 0.007079646017699115
$ws.Range("H12").Value = 0.03261882572227402
$ws.Range("D13").Value = 0
$ws.Range("F13").Value = 0.006489675516224189
$ws.Range("G13").Value = 0.01216182675601886
$ws.Range("I13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 0.008849557522123894
$ws.Range("I14").Value = 0.05003427004797804
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0.1045197740112993
$ws.Range("D16").Value = 0
$ws.Range("F16").Value = 0.1675516224188781
$ws.Range("G16").Value = 0.01389923057830726
$ws.Range("H16").Value = 0.2637465051258142
$ws.Range("I16").Value = 0
$ws.Range("F17").Value = 0.002949852507374631
$ws.Range("G17").Value = 0.02357905187391414
$ws.Range("I17").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0.01694915254237288
$ws.Range("D19").Value = 0
$ws.Range("F19").Value = 0.02713864306784661
$ws.Range("G19").Value = 0.006205013651030035
$ws.Range("I19").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("F20").Value = 0.01238938053097345
$ws.Range("G20").Value = 0.1092082402581289
$ws.Range("I20").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("F21").Value = 0.06371681415929203
$ws.Range("G21").Value = 0.1002730206006457
$ws.Range("I21").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0.258928571428571
$ws.Range("F23").Value = 0.01533923303834808
$ws.Range("G23").Value = 0.1015140233308517
$ws.Range("I23").Value = 0.1336531871144626
$ws.Range("B24").Value = 0.007949125596184421
$ws.Range("C24").Value = 0.02730696798493408
$ws.Range("F24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("F25").Value = 0.03303834808259587
$ws.Range("H25").Value = 0.06430568499534024
$ws.Range("D26").Value = 0
$ws.Range("F26").Value = 0.0005899705014749262
$ws.Range("D27").Value = 0
$ws.Range("F27").Value = 0.01061946902654867
$ws.Range("H27").Value = 0.03355079217148185
$ws.Range("D28").Value = 0
$ws.Range("F28").Value = 0.01769911504424779
$ws.Range("G28").Value = 0.02060064532141971
$ws.Range("I28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("F29").Value = 0.0359882005899705
$ws.Range("G29").Value = 0.000744601638123604
$ws.Range("I29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("F30").Value = 0.05722713864306785
$ws.Range("H30").Value = 0.01025163094128611
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0.05932203389830516
$ws.Range("H31").Value = 0.04659832246039147
$ws.Range("B32").Value = 0.6073131955484908
$ws.Range("C32").Value = 0.1468926553672317
$ws.Range("D32").Value = 0.6279069767441861
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("H32").Value = 0.0288909599254427
$ws.Range("B33").Value = 0.1049284578696345
$ws.Range("C33").Value = 0.1986817325800382
$ws.Range("D33").Value = 0
$ws.Range("F33").Value = 0.001769911504424779
$ws.Range("H33").Value = 0.09692451071761429
$ws.Range("G34").Value = 0.01067262347977165
$ws.Range("I34").Value = 0
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = 0.01506591337099811
$ws.Range("H35").Value = 0.04287045666356015